# Savesheet para diferentes bimestres
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value  = "-"
$ws.Range("D4").Value  = "MEC-1B-Gestão Integrada"
$ws.Range("E4").Value  = "-"
$ws.Range("D6").Value  = "MEC-1B-Gestão Integrada"
$ws.Range("C11").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("E14").Value = "MEC-1A-Gestão Integrada"
$ws.Range("E15").Value = "MEC-1A-Gestão Integrada"
